$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 2).Value = 36997
$ws.Cells.Item(3, 3).Value = 4954
$ws.Cells.Item(3, 4).Value = 6111
$ws.Cells.Item(4, 2).Value = 17127
$ws.Cells.Item(4, 3).Value = 1880
$ws.Cells.Item(4, 4).Value = 2149
$ws.Cells.Item(5, 2).Value = 59792
$ws.Cells.Item(5, 3).Value = 4154
$ws.Cells.Item(5, 4).Value = 5391
$ws.Cells.Item(6, 2).Value = 1337
$ws.Cells.Item(6, 3).Value = 386
$ws.Cells.Item(6, 4).Value = 35
$ws.Cells.Item(7, 2).Value = 35025
$ws.Cells.Item(7, 3).Value = 6276
$ws.Cells.Item(7, 4).Value = 5554
$ws.Cells.Item(8, 2).Value = 4075
$ws.Cells.Item(8, 3).Value = 902
$ws.Cells.Item(8, 4).Value = 851
$ws.Cells.Item(9, 2).Value = 4844
$ws.Cells.Item(9, 3).Value = 835
$ws.Cells.Item(9, 4).Value = 670
$ws.Cells.Item(10, 2).Value = 1750
$ws.Cells.Item(10, 3).Value = 241
$ws.Cells.Item(10, 4).Value = 58
$ws.Cells.Item(11, 2).Value = 345
$ws.Cells.Item(11, 3).Value = 186
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(13, 2).Value = 879
$ws.Cells.Item(13, 3).Value = 245
$ws.Cells.Item(13, 4).Value = 212
$ws.Cells.Item(14, 2).Value = 2353
$ws.Cells.Item(14, 3).Value = 1207
$ws.Cells.Item(14, 4).Value = 829
$ws.Cells.Item(15, 2).Value = 4779
$ws.Cells.Item(15, 3).Value = 1703
$ws.Cells.Item(15, 4).Value = 833
$ws.Cells.Item(16, 2).Value = 2933
$ws.Cells.Item(16, 3).Value = 1202
$ws.Cells.Item(16, 4).Value = 530
$ws.Cells.Item(17, 2).Value = 1564
$ws.Cells.Item(17, 3).Value = 551
$ws.Cells.Item(17, 4).Value = 111
$ws.Cells.Item(18, 2).Value = 14038
$ws.Cells.Item(18, 3).Value = 2132
$ws.Cells.Item(18, 4).Value = 2556
$ws.Cells.Item(19, 2).Value = 2481
$ws.Cells.Item(19, 3).Value = 638
$ws.Cells.Item(19, 4).Value = 448
$ws.Cells.Item(20, 2).Value = 17510
$ws.Cells.Item(20, 3).Value = 1825
$ws.Cells.Item(20, 4).Value = 3064
$ws.Cells.Item(21, 2).Value = 277
$ws.Cells.Item(21, 3).Value = 305
$ws.Cells.Item(21, 4).Value = 12
$ws.Cells.Item(22, 2).Value = 15261
$ws.Cells.Item(22, 3).Value = 1757
$ws.Cells.Item(22, 4).Value = 2340
$ws.Cells.Item(23, 2).Value = 949
$ws.Cells.Item(23, 3).Value = 391
$ws.Cells.Item(23, 4).Value = 140
$ws.Cells.Item(24, 2).Value = 14597
$ws.Cells.Item(24, 3).Value = 2581
$ws.Cells.Item(24, 4).Value = 2605
$ws.Cells.Item(25, 2).Value = 56062
$ws.Cells.Item(25, 3).Value = 6579
$ws.Cells.Item(25, 4).Value = 6348
$ws.Cells.Item(26, 2).Value = 4225
$ws.Cells.Item(26, 3).Value = 1309
$ws.Cells.Item(26, 4).Value = 554
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(28, 2).Value = 4785
$ws.Cells.Item(28, 3).Value = 894
$ws.Cells.Item(28, 4).Value = 991
$ws.Cells.Item(29, 2).Value = 2234
$ws.Cells.Item(29, 3).Value = 262
$ws.Cells.Item(29, 4).Value = 449
$ws.Cells.Item(30, 2).Value = 12903
$ws.Cells.Item(30, 3).Value = 2310
$ws.Cells.Item(30, 4).Value = 1946
$ws.Cells.Item(31, 2).Value = 423
$ws.Cells.Item(31, 3).Value = 223
$ws.Cells.Item(31, 4).Value = 173
$ws.Cells.Item(32, 2).Value = 1893
$ws.Cells.Item(32, 3).Value = 1325
$ws.Cells.Item(32, 4).Value = 229
$ws.Cells.Item(33, 2).Value = 12452
$ws.Cells.Item(33, 3).Value = 2564
$ws.Cells.Item(33, 4).Value = 2279
$ws.Cells.Item(34, 2).Value = 9547
$ws.Cells.Item(34, 3).Value = 2177
$ws.Cells.Item(34, 4).Value = 2399
$ws.Cells.Item(35, 2).Value = 5404
$ws.Cells.Item(35, 3).Value = 631
$ws.Cells.Item(35, 4).Value = 1185
$ws.Cells.Item(36, 2).Value = 42528
$ws.Cells.Item(36, 3).Value = 4554
$ws.Cells.Item(36, 4).Value = 4581
$ws.Cells.Item(37, 2).Value = 6978
$ws.Cells.Item(37, 3).Value = 2235
$ws.Cells.Item(37, 4).Value = 1105
$ws.Cells.Item(38, 2).Value = 15507
$ws.Cells.Item(38, 3).Value = 1647
$ws.Cells.Item(38, 4).Value = 2403
$ws.Cells.Item(39, 2).Value = 759
$ws.Cells.Item(39, 3).Value = 610
$ws.Cells.Item(39, 4).Value = 129
$ws.Cells.Item(40, 2).Value = 1351
$ws.Cells.Item(40, 3).Value = 210
$ws.Cells.Item(40, 4).Value = 753
$ws.Cells.Item(41, 2).Value = 1827
$ws.Cells.Item(41, 3).Value = 234
$ws.Cells.Item(41, 4).Value = 130
$ws.Cells.Item(42, 2).Value = 7719
$ws.Cells.Item(42, 3).Value = 248
$ws.Cells.Item(42, 4).Value = 192
$ws.Cells.Item(43, 2).Value = 240
$ws.Cells.Item(43, 3).Value = 102
$ws.Cells.Item(43, 4).Value = 72
$ws.Cells.Item(44, 2).Value = 551
$ws.Cells.Item(44, 3).Value = 48
$ws.Cells.Item(44, 4).Value = 22
$ws.Cells.Item(45, 2).Value = 1743
$ws.Cells.Item(45, 3).Value = 159
$ws.Cells.Item(45, 4).Value = 73
$ws.Cells.Item(46, 2).Value = 2729
$ws.Cells.Item(46, 3).Value = 768
$ws.Cells.Item(46, 4).Value = 401
$ws.Cells.Item(47, 2).Value = 9918
$ws.Cells.Item(47, 3).Value = 2619
$ws.Cells.Item(47, 4).Value = 2001
$ws.Cells.Item(48, 2).Value = 26083
$ws.Cells.Item(48, 3).Value = 2566
$ws.Cells.Item(48, 4).Value = 3867
$ws.Cells.Item(49, 2).Value = 12412
$ws.Cells.Item(49, 3).Value = 2710
$ws.Cells.Item(49, 4).Value = 902
$ws.Cells.Item(50, 2).Value = 8730
$ws.Cells.Item(50, 3).Value = 875
$ws.Cells.Item(50, 4).Value = 1418
$ws.Cells.Item(51, 2).Value = 21707
$ws.Cells.Item(51, 3).Value = 2443
$ws.Cells.Item(51, 4).Value = 2637
$ws.Cells.Item(52, 2).Value = 3183
$ws.Cells.Item(52, 3).Value = 380
$ws.Cells.Item(52, 4).Value = 689
$ws.Cells.Item(53, 2).Value = 11102
$ws.Cells.Item(53, 3).Value = 2274
$ws.Cells.Item(53, 4).Value = 2021
$ws.Cells.Item(54, 2).Value = 1882
$ws.Cells.Item(54, 3).Value = 565
$ws.Cells.Item(54, 4).Value = 972
$ws.Cells.Item(55, 2).Value = 1580
$ws.Cells.Item(55, 3).Value = 1131
$ws.Cells.Item(55, 4).Value = 192
$ws.Cells.Item(56, 2).Value = 2978
$ws.Cells.Item(56, 3).Value = 864
$ws.Cells.Item(56, 4).Value = 1037
$ws.Cells.Item(57, 2).Value = 10944
$ws.Cells.Item(57, 3).Value = 4357
$ws.Cells.Item(57, 4).Value = 2404
$ws.Cells.Item(58, 2).Value = 12080
$ws.Cells.Item(58, 3).Value = 791
$ws.Cells.Item(58, 4).Value = 489
$ws.Cells.Item(59, 2).Value = 530348
$ws.Cells.Item(59, 3).Value = 79867
$ws.Cells.Item(59, 4).Value = 76650
